$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for 2025-05-15 / SAT BRAS INDUSTRIA ELETRONICA DA AMAZONIA LTDA. / FITA DUREX (old row 2)
$ws.Rows.Item(2).Delete()

# The "A" column holds a pre-computed index value (0-based) that, in the source
# dataset, is recomputed after the row removal: every value greater than the
# removed row's original index (7) is decremented by one.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -gt 7) {
        $cell.Value = $val - 1
    }
}
